# Adicao de algumas funcoes adicionais reponsaveis pela conversao do formato das rotas.
# This updates columns E (distance/cost) and F (ratio) for rows 1-15,
# extending the data from rows 1-8 down through rows 9-15 (new rows),
# and growing the worksheet's used range / dimension to A1:F15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = 3918.965007467933
$ws.Range("F1").Value = 0.8744726293774652

$ws.Range("E2").Value = 4046.226505918316
$ws.Range("F2").Value = 0.8771926484679836

$ws.Range("E3").Value = 4000.938865799878
$ws.Range("F3").Value = 0.8761829498419887

$ws.Range("E4").Value = 3657.135991387497
$ws.Range("F4").Value = 0.8638507711238541

$ws.Range("E5").Value = 3691.662049029868
$ws.Range("F5").Value = 0.8667998525276273

$ws.Range("E6").Value = 3743.48109539564
$ws.Range("F6").Value = 0.8710824545673868

$ws.Range("E7").Value = 4075.612340139283
$ws.Range("F7").Value = 0.8773679444291423

$ws.Range("E8").Value = 3649.317874322935
$ws.Range("F8").Value = 0.8609585078099218

$ws.Range("E9").Value = 3936.104497130708
$ws.Range("F9").Value = 0.8750245789681041

$ws.Range("E10").Value = 3957.666379190404
$ws.Range("F10").Value = 0.875905878642028

$ws.Range("E11").Value = 3773.551720558038
$ws.Range("F11").Value = 0.8718342842150564

$ws.Range("E12").Value = 4002.437914321994
$ws.Range("F12").Value = 0.8765861657046756

$ws.Range("E13").Value = 3884.662082949515
$ws.Range("F13").Value = 0.8743602821837185

$ws.Range("E14").Value = 3939.123864467525
$ws.Range("F14").Value = 0.8758262600701566

$ws.Range("E15").Value = 3794.989528581777
$ws.Range("F15").Value = 0.8737554517519703
